# Fruta / hortaliza, semanal
#
# A new weekly observation is inserted as row 18 (pushing the existing
# rows 18-57 down to 19-58). The new row carries the same market /
# category metadata as its neighbours, with its own date, volume,
# price and origin figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 18.. down by one to make room for the new observation.
$ws.Range("A18").EntireRow.Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44526
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112026
$ws.Range("G18").Value = "Haba"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 8000
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Región de La Araucanía"
$ws.Range("P18").Value = 320
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
